$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name
$ws.Range("C3").Value = "Beerdavinder singh"

# Row 7: __init__ / attributes set to input values
$ws.Range("E7").Value = "Inputs"
$ws.Range("F7").Value = "700, 10, 100.00, date(2024, 10, 14), -100.00, 0.05"
$ws.Range("G7").Value = "Attributes setup"

# Row 8: overdraft limit invalid type
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = """bs"" (instead of -100.00)"
$ws.Range("G8").Value = "Raises ValueError"

# Row 9: overdraft rate invalid type
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = """bs"" (instead of 0.05)"
$ws.Range("G9").Value = "Raises ValueError"

# Row 10: date created invalid type
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = """bs"" (instead of date(2024, 10, 14))"
$ws.Range("G10").Value = "Raises ValueError"

# Row 11: balance greater than overdraft limit
$ws.Range("E11").Value = "Balance = `$100.00, Overdraft limit = `$-100.00"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "Returns 0.50"

# Row 12: balance less than overdraft limit
$ws.Range("E12").Value = "Balance = `$-200.00, Overdraft limit = `$-100.00"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "`tReturns 5.50"

# Row 13: balance equal to overdraft limit
$ws.Range("E13").Value = "Balance = `$-100.00, Overdraft limit = `$-100.00"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "Returns 0.50"

# Row 14: __str__ / appropriate value returned
$ws.Range("E14").Value = "700, 10, 100.00, date(2024, 10, 14), -100.00, 0.05"
$ws.Range("F14").Value = "700, 10, 100.00, date(2024, 10, 14), -100.00, 0.05"
$ws.Range("G14").Value = "Format str"

$ws.Range("F13").Select()
